$wb = $excel.ActiveWorkbook

$labels = $wb.Worksheets.Item("Labels")

$labels.Range("I9").Value = "PF :"
$labels.Range("I8").Value = "Stocks : "
$labels.Range("I7").Value = "Mutual Funds :"
$labels.Range("I6").Value = "Fixed Deposits :"
$labels.Range("I10").Value = "Crypto-Currency :"
$labels.Range("I5").Value = "Date :"
$labels.Range("I4").Value = "Month :"
$labels.Range("I3").Value = "Year :"

$labels.Activate()
$labels.Range("A10").Select()
